$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel stores them as text (matching the inlineStr type in the source data)
# rather than auto-converting to a numeric cell type.
$ws.Range('D2').Value = '64.980.13'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '3.163.48'
$ws.Range('E3').Value = '  +3.60%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.71'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.48'
$ws.Range('E6').Value = '  +5.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.162.54'
$ws.Range('E8').Value = '  +3.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  +3.07%  '
$ws.Range('E10').Value = '  +4.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.20'
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.505'
$ws.Range('E12').Value = '  +5.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000273'
$ws.Range('E13').Value = '  +18.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.15'
$ws.Range('E14').Value = '  +8.22%  '
$ws.Range('D15').Value = '3.678.01'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('D16').Value = '65.049.66'
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.161.41'
$ws.Range('E17').Value = '  +3.49%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.18'
$ws.Range('E18').Value = '  +6.49%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '510.02'
$ws.Range('E20').Value = '  +6.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.91'
$ws.Range('E21').Value = '  +6.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.72'
$ws.Range('E22').Value = '  +9.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.733'
$ws.Range('E23').Value = '  +8.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.82'
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.73'
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +15.24%  '
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  +8.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.09'
$ws.Range('E30').Value = '  +6.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.79'
$ws.Range('E31').Value = '  +15.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +7.33%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.33'
$ws.Range('E34').Value = '  +11.98%  '
$ws.Range('E35').Value = '  +7.08%  '
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '476.03'
$ws.Range('E37').Value = '  +7.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0880'
$ws.Range('E38').Value = '  +9.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.08'
$ws.Range('E39').Value = '  +9.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0419'
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('D41').Value = '3.120.05'
$ws.Range('E41').Value = '  +4.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.62'
$ws.Range('E42').Value = '  +4.69%  '
$ws.Range('E43').Value = '  +5.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  +17.90%  '
$ws.Range('E45').Value = '  +10.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.31'
$ws.Range('E46').Value = '  +5.45%  '
$ws.Range('D47').Value = '0.0₃0586'
$ws.Range('E47').Value = '  +13.40%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.115'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('E50').Value = '  +11.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '122.84'
$ws.Range('E51').Value = '  +3.89%  '

Write-Output "Applied 92 cell updates"